$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.89
$ws.Range("J2").Value = 2.9
$ws.Range("Q2").Value = 1.86
$ws.Range("W2").Value = 2.12
$ws.Range("X2").Value = 16.5
$ws.Range("AB2").Value = 8.800000000000001
$ws.Range("AE2").Value = 1000
$ws.Range("AK2").Value = 21

# Row 3
$ws.Range("F3").Value = 1.98
$ws.Range("L3").Value = 1.28
$ws.Range("AA3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AK3").Value = 24
$ws.Range("AL3").Value = 34
$ws.Range("AM3").Value = 1000

# Row 4
$ws.Range("K4").Value = 3.75
$ws.Range("W4").Value = 1.51
$ws.Range("AB4").Value = 1000

# Row 6
$ws.Range("F6").Value = 3.55
$ws.Range("G6").Value = 4.8
$ws.Range("H6").Value = 2.12
$ws.Range("I6").Value = 2.56
$ws.Range("J6").Value = 2.98
$ws.Range("K6").Value = 3.95
$ws.Range("M6").Value = 1.09
$ws.Range("N6").Value = 2.56
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 1.62
$ws.Range("Q6").Value = 2.12
$ws.Range("R6").Value = 1.23
$ws.Range("S6").Value = 4.7
$ws.Range("T6").Value = 1.94
$ws.Range("U6").Value = 1.84
$ws.Range("V6").Value = 1.64
$ws.Range("W6").Value = 1.29
